$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7768072530323593
$ws.Range("C2").Value = 0.1947694543560772
$ws.Range("E2").Value = 0.7294059226006908
$ws.Range("F2").Value = 2.125094005473898
$ws.Range("G2").Value = 0.002389931180061565
$ws.Range("I2").Value = 0.3279309542900108
$ws.Range("J2").Value = 0.03163005799406804
$ws.Range("N2").Value = 0.8905864980084175
$ws.Range("O2").Value = 1.62647748621194
$ws.Range("B3").Value = 0.6840290669999263
$ws.Range("C3").Value = 0.1701091766746856
$ws.Range("E3").Value = 0.7065446905868527
$ws.Range("F3").Value = 2.087326923490849
$ws.Range("G3").Value = 0.002392582710328061
$ws.Range("I3").Value = 0.3330822552223189
$ws.Range("J3").Value = 0.0316440336508883
$ws.Range("N3").Value = 0.8893806391606063
$ws.Range("O3").Value = 1.626116702714739
$ws.Range("B4").Value = 0.6269627375471316
$ws.Range("C4").Value = 0.1549056167598906
$ws.Range("E4").Value = 0.6928565713068195
$ws.Range("F4").Value = 2.065492822590727
$ws.Range("G4").Value = 0.002394296713355615
$ws.Range("I4").Value = 0.3365238082549205
$ws.Range("J4").Value = 0.03168008694556335
$ws.Range("N4").Value = 0.8890156574449719
$ws.Range("O4").Value = 1.627291173391029
$ws.Range("B5").Value = 0.6036839437693118
$ws.Range("C5").Value = 0.1486947258585758
$ws.Range("E5").Value = 0.6873661571158607
$ws.Range("F5").Value = 2.056935217299582
$ws.Range("G5").Value = 0.002395016863060957
$ws.Range("I5").Value = 0.3379961413367489
$ws.Range("J5").Value = 0.03170166262354357
$ws.Range("N5").Value = 0.8889615947955676
$ws.Range("O5").Value = 1.628119772097534
$ws.Range("B6").Value = 0.5998171210544001
$ws.Range("C6").Value = 0.147662495517892
$ws.Range("E6").Value = 0.6864597664936269
$ws.Range("F6").Value = 2.055534744422786
$ws.Range("G6").Value = 0.002395137754742076
$ws.Range("I6").Value = 0.33824483560376
$ws.Range("J6").Value = 0.0317056602915109
$ws.Range("N6").Value = 0.8889583435798656
$ws.Range("O6").Value = 1.628278468504263
$ws.Range("B7").Value = 0.6266488861288622
$ws.Range("C7").Value = 0.1548219161487339
$ws.Range("E7").Value = 0.6927821710743558
$ws.Range("F7").Value = 2.065376036236074
$ws.Range("G7").Value = 0.002394306337752373
$ws.Range("I7").Value = 0.336543381999693
$ws.Range("J7").Value = 0.03168035008279801
$ws.Range("N7").Value = 0.8890145446629987
$ws.Range("O7").Value = 1.627300932443404
$ws.Range("B8").Value = 0.7448389195902791
$ws.Range("C8").Value = 0.1862796330416359
$ws.Range("E8").Value = 0.7214509157749092
$ws.Range("F8").Value = 2.111790195963636
$ws.Range("G8").Value = 0.002390827627995393
$ws.Range("I8").Value = 0.329649177111861
$ws.Range("J8").Value = 0.0316291601310752
$ws.Range("N8").Value = 0.8900929557204051
$ws.Range("O8").Value = 1.626062717091884
$ws.Range("B9").Value = 0.9757651099570808
$ws.Range("C9").Value = 0.2474661312647868
$ws.Range("E9").Value = 0.7804456863386804
$ws.Range("F9").Value = 2.213602304440997
$ws.Range("G9").Value = 0.00238468485700918
$ws.Range("I9").Value = 0.3183490719464714
$ws.Range("J9").Value = 0.0317479082976746
$ws.Range("N9").Value = 0.8951760698139992
$ws.Range("O9").Value = 1.634760864569216
$ws.Range("B10").Value = 1.144862061827439
$ws.Range("C10").Value = 0.2921051600833096
$ws.Range("E10").Value = 0.8254978880137003
$ws.Range("F10").Value = 2.29505406678976
$ws.Range("G10").Value = 0.002380581429634146
$ws.Range("I10").Value = 0.3114116267707701
$ws.Range("J10").Value = 0.03197040741874346
$ws.Range("N10").Value = 0.9007081368950338
$ws.Range("O10").Value = 1.648006998776566
$ws.Range("B11").Value = 1.221656550196769
$ws.Range("C11").Value = 0.3123427204654092
$ws.Range("E11").Value = 0.846368334989748
$ws.Range("F11").Value = 2.333568533435567
$ws.Range("G11").Value = 0.0023788027239981
$ws.Range("I11").Value = 0.3085544677278165
$ws.Range("J11").Value = 0.03210136144416964
$ws.Range("N11").Value = 0.9036127953959578
$ws.Range("O11").Value = 1.655537553615233
$ws.Range("B12").Value = 1.250716970956319
$ws.Range("C12").Value = 0.3199960141935492
$ws.Range("E12").Value = 0.8543257080803102
$ws.Range("F12").Value = 2.348364210146457
$ws.Range("G12").Value = 0.002378141755295467
$ws.Range("I12").Value = 0.3075157031674003
$ws.Range("J12").Value = 0.03215525494960758
$ws.Range("N12").Value = 0.9047683073592054
$ws.Range("O12").Value = 1.65860678880415
$ws.Range("B13").Value = 1.244459202027144
$ws.Range("C13").Value = 0.3183482004891118
$ws.Range("E13").Value = 0.8526095333919841
$ws.Range("F13").Value = 2.345168293994817
$ws.Range("G13").Value = 0.002378283547937661
$ws.Range("I13").Value = 0.3077374957090875
$ws.Range("J13").Value = 0.03214345615200997
$ws.Range("N13").Value = 0.904516978752099
$ws.Range("O13").Value = 1.657936078375855
$ws.Range("B14").Value = 1.224047777596752
$ws.Range("C14").Value = 0.3129725683934623
$ws.Range("E14").Value = 0.8470219061403839
$ws.Range("F14").Value = 2.334781547039512
$ws.Range("G14").Value = 0.002378748093827416
$ws.Range("I14").Value = 0.3084681412756538
$ws.Range("J14").Value = 0.03210570887639008
$ws.Range("N14").Value = 0.9037067474141338
$ws.Range("O14").Value = 1.655785694067902
$ws.Range("B15").Value = 1.211542543730843
$ws.Range("C15").Value = 0.3096784967901556
$ws.Range("E15").Value = 0.8436063836765868
$ws.Range("F15").Value = 2.328446880848276
$ws.Range("G15").Value = 0.002379034279086256
$ws.Range("I15").Value = 0.3089213126878292
$ws.Range("J15").Value = 0.03208314892074426
$ws.Range("N15").Value = 0.9032176889707841
$ws.Range("O15").Value = 1.654496891815285
$ws.Range("B16").Value = 1.139840646183643
$ws.Range("C16").Value = 0.2907811733388712
$ws.Range("E16").Value = 0.8241415312233045
$ws.Range("F16").Value = 2.292566544770153
$ws.Range("G16").Value = 0.002380699436958995
$ws.Range("I16").Value = 0.3116043789837164
$ws.Range("J16").Value = 0.03196245004289011
$ws.Range("N16").Value = 0.9005261005563625
$ws.Range("O16").Value = 1.647545241528007
$ws.Range("B17").Value = 1.095819858443804
$ws.Range("C17").Value = 0.2791703910248202
$ws.Range("E17").Value = 0.812296857161158
$ws.Range("F17").Value = 2.270930122960152
$ws.Range("G17").Value = 0.002381743443620595
$ws.Range("I17").Value = 0.3133270218654118
$ws.Range("J17").Value = 0.03189604155264902
$ws.Range("N17").Value = 0.8989741351347789
$ws.Range("O17").Value = 1.643666893777208
$ws.Range("B18").Value = 1.070488253817928
$ws.Range("C18").Value = 0.2724857042475151
$ws.Range("E18").Value = 0.805519494821624
$ws.Range("F18").Value = 2.258622948623838
$ws.Range("G18").Value = 0.002382352211998649
$ws.Range("I18").Value = 0.3143459439790988
$ws.Range("J18").Value = 0.03186064256943055
$ws.Range("N18").Value = 0.8981180280017185
$ws.Range("O18").Value = 1.641577740715547
$ws.Range("B19").Value = 1.061909390458311
$ws.Range("C19").Value = 0.2702212822376282
$ws.Range("E19").Value = 0.803230869401375
$ws.Range("F19").Value = 2.254479543308264
$ws.Range("G19").Value = 0.002382559754875252
$ws.Range("I19").Value = 0.3146957534785919
$ws.Range("J19").Value = 0.03184913669742784
$ws.Range("N19").Value = 0.897834448023886
$ws.Range("O19").Value = 1.64089466805666
$ws.Range("B20").Value = 1.100507200497304
$ws.Range("C20").Value = 0.2804070507278311
$ws.Range("E20").Value = 0.8135540802659307
$ws.Range("F20").Value = 2.273219117356177
$ws.Range("G20").Value = 0.00238163145045358
$ws.Range("I20").Value = 0.3131407333602922
$ws.Range("J20").Value = 0.03190282112625198
$ws.Range("N20").Value = 0.8991355641677785
$ws.Range("O20").Value = 1.644065089959327
$ws.Range("B21").Value = 1.230043662036223
$ws.Range("C21").Value = 0.3145518020211284
$ws.Range("E21").Value = 0.8486616560693818
$ws.Range("F21").Value = 2.337826650423267
$ws.Range("G21").Value = 0.002378611304292691
$ws.Range("I21").Value = 0.308252359294297
$ws.Range("J21").Value = 0.03211667913414828
$ws.Range("N21").Value = 0.9039432253252215
$ws.Range("O21").Value = 1.656411399617411
$ws.Range("B22").Value = 1.314586047322734
$ws.Range("C22").Value = 0.3368075906011541
$ws.Range("E22").Value = 0.8719223940838816
$ws.Range("F22").Value = 2.381282191041208
$ws.Range("G22").Value = 0.002376710811024182
$ws.Range("I22").Value = 0.3053093118736925
$ws.Range("J22").Value = 0.03228154798041061
$ws.Range("N22").Value = 0.9074091766364205
$ws.Range("O22").Value = 1.665749185988403
$ws.Range("B23").Value = 1.269475400610418
$ws.Range("C23").Value = 0.3249348309883828
$ws.Range("E23").Value = 0.8594787531639696
$ws.Range("F23").Value = 2.357976234938775
$ws.Range("G23").Value = 0.00237771844869302
$ws.Range("I23").Value = 0.3068569590162618
$ws.Range("J23").Value = 0.03219124850687649
$ws.Range("N23").Value = 0.905529772615111
$ws.Range("O23").Value = 1.660648950008891
$ws.Range("B24").Value = 1.098388126138389
$ws.Range("C24").Value = 0.2798479864259207
$ws.Range("E24").Value = 0.8129855891032491
$ws.Range("F24").Value = 2.272183852211185
$ws.Range("G24").Value = 0.002381682056004994
$ws.Range("I24").Value = 0.3132248654714083
$ws.Range("J24").Value = 0.03189974742387136
$ws.Range("N24").Value = 0.8990624695124438
$ws.Range("O24").Value = 1.64388462770745
$ws.Range("B25").Value = 0.9133890096939012
$ws.Range("C25").Value = 0.2309682053049471
$ws.Range("E25").Value = 0.7641871879279591
$ws.Range("F25").Value = 2.184896637946764
$ws.Range("G25").Value = 0.002386274392122229
$ws.Range("I25").Value = 0.3211672351900248
$ws.Range("J25").Value = 0.03169215225198485
$ws.Range("N25").Value = 0.8934844439991991
$ws.Range("O25").Value = 1.631209127962109
